$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'67.792.34"
$ws.Range("E2").Value = "  +0.27%  "
# Row 3
$ws.Range("D3").Value = "'3.812.93"
$ws.Range("E3").Value = "  +0.71%  "
# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.08%  "
# Row 5
$ws.Range("D5").Value = "'608.14"
$ws.Range("E5").Value = "  +2.11%  "
# Row 6
$ws.Range("D6").Value = "'167.01"
$ws.Range("E6").Value = "  +0.21%  "
# Row 7
$ws.Range("E7").Value = "  -0.01%  "
# Row 8
$ws.Range("E8").Value = "  +0.51%  "
# Row 9
$ws.Range("E9").Value = "  +0.70%  "
# Row 10
$ws.Range("E10").Value = "  -0.76%  "
# Row 11
$ws.Range("E11").Value = "  +0.57%  "
# Row 12
$ws.Range("D12").Value = "'0.0000254"
$ws.Range("E12").Value = "  -0.90%  "
# Row 13
$ws.Range("D13").Value = "'36.09"
$ws.Range("E13").Value = "  -0.78%  "
# Row 14
$ws.Range("D14").Value = "'4.450.31"
$ws.Range("E14").Value = "  +0.49%  "
# Row 15
$ws.Range("D15").Value = "'3.829.01"
$ws.Range("E15").Value = "  +1.36%  "
# Row 16
$ws.Range("D16").Value = "'18.52"
$ws.Range("E16").Value = "  -0.06%  "
# Row 17
$ws.Range("D17").Value = "'67.790.77"
$ws.Range("E17").Value = "  +0.29%  "
# Row 18
$ws.Range("D18").Value = "'7.10"
$ws.Range("E18").Value = "  +1.45%  "
# Row 19
$ws.Range("E19").Value = "  +0.48%  "
# Row 20
$ws.Range("D20").Value = "'462.55"
$ws.Range("E20").Value = "  +1.15%  "
# Row 21
$ws.Range("D21").Value = "'9.89"
$ws.Range("E21").Value = "  -2.97%  "
# Row 22
$ws.Range("D22").Value = "'0.702"
$ws.Range("E22").Value = "  +0.65%  "
# Row 23
$ws.Range("E23").Value = "  -2.99%  "
# Row 24
$ws.Range("D24").Value = "'83.44"
$ws.Range("E24").Value = "  -0.09%  "
# Row 25
$ws.Range("D25").Value = "'12.10"
$ws.Range("E25").Value = "  +1.64%  "
# Row 26
$ws.Range("E26").Value = "  -1.25%  "
# Row 27
$ws.Range("E27").Value = "  +0.21%  "
# Row 28
$ws.Range("D28").Value = "'10.01"
$ws.Range("E28").Value = "  -0.54%  "
# Row 29
$ws.Range("D29").Value = "'3.962.21"
$ws.Range("E29").Value = "  +0.55%  "
# Row 30
$ws.Range("E30").Value = "  +0.61%  "
# Row 31
$ws.Range("D31").Value = "'7.41"
$ws.Range("E31").Value = "  +1.65%  "
# Row 32
$ws.Range("D32").Value = "'2.24"
$ws.Range("E32").Value = "  +1.62%  "
# Row 33
$ws.Range("D33").Value = "'29.62"
$ws.Range("E33").Value = "  -0.73%  "
# Row 34
$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").Value = "'9.08"
$ws.Range("E34").Value = "  -1.55%  "
# Row 35
$ws.Range("B35").Value = "RenzoRestakedETH"
$ws.Range("C35").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D35").Value = "'3.757.61"
$ws.Range("E35").Value = "  +0.44%  "
# Row 36
$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D36").Value = "'0.993"
$ws.Range("E36").Value = "  -0.50%  "
# Row 37
$ws.Range("E37").Value = "  -0.01%  "
# Row 38
$ws.Range("E38").Value = "  +1.54%  "
# Row 39
$ws.Range("E39").Value = "  -0.06%  "
# Row 40
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.37%  "
# Row 41
$ws.Range("D41").Value = "'5.80"
$ws.Range("E41").Value = "  +0.72%  "
# Row 42
$ws.Range("E42").Value = "  -0.13%  "
# Row 44
$ws.Range("D44").Value = "'48.14"
$ws.Range("E44").Value = "  +2.36%  "
# Row 45
$ws.Range("E45").Value = "  +0.60%  "
# Row 46
$ws.Range("D46").Value = "'43.14"
$ws.Range("E46").Value = "  -4.17%  "
# Row 47
$ws.Range("D47").Value = "'28.01"
$ws.Range("E47").Value = "  +9.32%  "
# Row 48
$ws.Range("D48").Value = "'8.36"
$ws.Range("E48").Value = "  +0.06%  "
# Row 49
$ws.Range("D49").Value = "'148.85"
$ws.Range("E49").Value = "  -0.04%  "
# Row 50
$ws.Range("E50").Value = "  +10.00%  "
# Row 51
$ws.Range("E51").Value = "  +0.25%  "
